$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.877.55"
$ws.Range("E2").Value = "  +2.22%  "
$ws.Range("D3").Value = "3.094.78"
$ws.Range("E3").Value = "  +4.93%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'579.13"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("D6").Value = "'170.27"
$ws.Range("E6").Value = "  +5.29%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.089.97"
$ws.Range("E8").Value = "  +4.82%  "
$ws.Range("D9").Value = "'0.523"
$ws.Range("E9").Value = "  +1.44%  "
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("E11").Value = "  +2.90%  "
$ws.Range("D12").Value = "'0.479"
$ws.Range("E12").Value = "  +4.73%  "
$ws.Range("D13").Value = "'0.0000249"
$ws.Range("E13").Value = "  +1.86%  "
$ws.Range("D14").Value = "'36.64"
$ws.Range("E14").Value = "  +6.22%  "
$ws.Range("E15").Value = "  -0.76%  "
$ws.Range("D16").Value = "3.607.06"
$ws.Range("E16").Value = "  +4.90%  "
$ws.Range("D17").Value = "66.845.26"
$ws.Range("E17").Value = "  +2.22%  "
$ws.Range("D18").Value = "'7.18"
$ws.Range("E18").Value = "  +1.62%  "
$ws.Range("D19").Value = "3.095.60"
$ws.Range("E19").Value = "  +4.92%  "
$ws.Range("D20").Value = "'16.21"
$ws.Range("E20").Value = "  +4.36%  "
$ws.Range("D21").Value = "'465.76"
$ws.Range("E21").Value = "  +4.73%  "
$ws.Range("E22").Value = "  +2.83%  "
$ws.Range("D23").Value = "'7.47"
$ws.Range("E23").Value = "  +2.25%  "
$ws.Range("D24").Value = "'83.86"
$ws.Range("E24").Value = "  +1.52%  "
$ws.Range("D25").Value = "'13.15"
$ws.Range("E25").Value = "  +8.03%  "
$ws.Range("D26").Value = "'2.33"
$ws.Range("E26").Value = "  +2.29%  "
$ws.Range("D27").Value = "'10.06"
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").Value = "'8.00"
$ws.Range("E29").Value = "  -1.20%  "
$ws.Range("D30").Value = "'2.40"
$ws.Range("E30").Value = "  -2.30%  "
$ws.Range("E31").Value = "  +3.29%  "
$ws.Range("D32").Value = "'0.0000104"
$ws.Range("E32").Value = "  +2.76%  "
$ws.Range("D33").Value = "'28.23"
$ws.Range("E33").Value = "  +3.68%  "
$ws.Range("E34").Value = "  +1.64%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +2.80%  "
$ws.Range("E37").Value = "  +2.61%  "
$ws.Range("D38").Value = "'47.64"
$ws.Range("E38").Value = "  +4.95%  "
$ws.Range("D39").Value = "'2.11"
$ws.Range("E39").Value = "  +8.76%  "
$ws.Range("D40").Value = "'50.22"
$ws.Range("E40").Value = "  +1.49%  "
$ws.Range("D41").Value = "'0.316"
$ws.Range("E41").Value = "  +4.62%  "
$ws.Range("E42").Value = "  +1.34%  "
$ws.Range("D43").Value = "'8.65"
$ws.Range("E43").Value = "  +0.95%  "
$ws.Range("D44").Value = "'2.80"
$ws.Range("E44").Value = "  -1.55%  "
$ws.Range("D45").Value = "'0.0359"
$ws.Range("E45").Value = "  +1.98%  "
$ws.Range("D46").Value = "'381.88"
$ws.Range("E46").Value = "  -0.65%  "
$ws.Range("D47").Value = "2.790.50"
$ws.Range("E47").Value = "  +3.42%  "
$ws.Range("D48").Value = "'134.57"
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("D50").Value = "'24.69"
$ws.Range("E50").Value = "  +4.85%  "
$ws.Range("D51").Value = "'2.23"
$ws.Range("E51").Value = "  +1.82%  "
